$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AN1").Value = "Switzerland vs Italy"
